$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the reported accuracy figures (B15, B18) to match the latest run.
$ws.Range("B15").Value = "Accuracy on test set : 69.0037%"
$ws.Range("B18").Value = "`nprediction accuracy on test set: 68.2657%`n"

# Move the saved selection from B19 to A18, as in the uploaded workbook.
$ws.Range("A18").Select() | Out-Null
